$d = $word.ActiveDocument

# --- "Programa resumido" / "Programa" sections (Portuguese) ---
# Same text appears twice in the document; ReplaceAll handles both.
$oldPt = "1. Problemas ambientais, causas e soluções2. Ecossistemas: o que são e como funcionam3. Biodiversidade e evolução4. Biodiversidade, interações de espécies e controle da população5. A população humana e seu impacto"
$newPt = "1. Problemas ambientais, causas e soluções^l2. Ecossistemas: o que são e como funcionam^l3. Biodiversidade e evolução^l4. Biodiversidade, interações de espécies e controle da população^l5. A população humana e seu impacto"
$d.Content.Find.Execute($oldPt, $true, $false, $false, $false, $false, $true, 1, $false, $newPt, 2)

# --- "Programa resumido" / "Programa" sections (English, italic) ---
$oldEn = "1. Environmental problems, causes and solutions2. Ecosystems: what they are and how they work3. Biodiversity and evolution4. Biodiversity, species interactions and population control5. The human population and its impact"
$newEn = "1. Environmental problems, causes and solutions^l2. Ecosystems: what they are and how they work^l3. Biodiversity and evolution^l4. Biodiversity, species interactions and population control^l5. The human population and its impact"
$d.Content.Find.Execute($oldEn, $true, $false, $false, $false, $false, $true, 1, $false, $newEn, 2)

# --- "Avaliação" -> "Critério" run ---
$oldCriterio = "O aluno poderá optar por dois critérios de avaliação:Critério 1: NF = (P1+P2)/2; ouCritério 2: NF = (NOTA 1 + NOTA 2)/2Sendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$newCriterio = "O aluno poderá optar por dois critérios de avaliação:^lCritério 1: NF = (P1+P2)/2; ou^lCritério 2: NF = (NOTA 1 + NOTA 2)/2^lSendo P1 e P2 avaliações escritas e NOTA 1 e NOTA 2 obtidas em atividades desenvolvidas em aula, trabalhos e relatórios de aulas práticas."
$d.Content.Find.Execute($oldCriterio, $true, $false, $false, $false, $false, $true, 1, $false, $newCriterio, 2)

# --- "Bibliografia" paragraph ---
$oldBib = "Básica:MILLER, G.T.; SPOOLMAN, S.E. 2012. Ecologia e sustentabilidade. Cengage Learning. 412p.Complementar:BEGON, M., J.L. HARPER & C.R. TOWNSEND. 2005. Ecology. From Individuals to Communities. Blackwell Science.RICKLEFS, R.E. 2003. A economia da natureza. Guanabara Koogan.RICKLEFS, R.E. & G.L. MILLER. 2000. Ecology. W.H. Freeman and Co.TOWNSEND, C.R., M. BEGON. & J.L. HARPER 2006. Fundamentos em ecologia. Artmed."
$newBib = "Básica:^lMILLER, G.T.; SPOOLMAN, S.E. 2012. Ecologia e sustentabilidade. Cengage Learning. 412p.^l^lComplementar:^lBEGON, M., J.L. HARPER & C.R. TOWNSEND. 2005. Ecology. From Individuals to Communities. Blackwell Science.^lRICKLEFS, R.E. 2003. A economia da natureza. Guanabara Koogan.^lRICKLEFS, R.E. & G.L. MILLER. 2000. Ecology. W.H. Freeman and Co.^lTOWNSEND, C.R., M. BEGON. & J.L. HARPER 2006. Fundamentos em ecologia. Artmed."
$d.Content.Find.Execute($oldBib, $true, $false, $false, $false, $false, $true, 1, $false, $newBib, 2)
